# Update "想去人数" (interested-count) figures that changed between crawls.
# Same five rows are refreshed on both the "展览" and "全部类型" sheets,
# which mirror each other in this workbook.

$wb = $excel.ActiveWorkbook

$updates = @{
    8  = 11295
    11 = 26
    13 = 2512
    19 = 11257
    20 = 11114
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
